$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 54.114
$ws.Range("D2").Value = 54.114
$ws.Range("E2").Value = 3.16948242
$ws.Range("F2").Value = 0.00135841
$ws.Range("G2").Value = 0.07295569
$ws.Range("H2").Value = 3.97650341
$ws.Range("I2").Value = 6.551482727285933
$ws.Range("J2").Value = 6.551482727285933
$ws.Range("K2").Value = 0.3907925152580654
$ws.Range("L2").Value = 0.0001665157632159772
$ws.Range("M2").Value = 0.008893789057972293
$ws.Range("N2").Value = 0.8448462905931474
$ws.Range("C3").Value = 90.295
$ws.Range("D3").Value = 90.295
$ws.Range("E3").Value = 1.904231
$ws.Range("F3").Value = 0.000943
$ws.Range("G3").Value = 0.08433705000000001
$ws.Range("H3").Value = 7.67577954
$ws.Range("I3").Value = 11.89283037861107
$ws.Range("J3").Value = 11.89283037861107
$ws.Range("K3").Value = 0.2527369708595069
$ws.Range("L3").Value = 0.0001822350394671936
$ws.Range("M3").Value = 0.01559989991233196
$ws.Range("N3").Value = 2.005840126932154
$ws.Range("C4").Value = 27.812
$ws.Range("D4").Value = 55.605
$ws.Range("E4").Value = 3.104128
$ws.Range("F4").Value = 0.00264691
$ws.Range("G4").Value = 0.0365285
$ws.Range("H4").Value = 1.03036605
$ws.Range("I4").Value = 4.142924066383527
$ws.Range("J4").Value = 8.279222217124053
$ws.Range("K4").Value = 0.4427500415151439
$ws.Range("L4").Value = 0.0005123196495580892
$ws.Range("M4").Value = 0.007658787474240003
$ws.Range("N4").Value = 0.3245784489303126
$ws.Range("C5").Value = 46.208
$ws.Range("D5").Value = 90.684
$ws.Range("E5").Value = 1.91382706
$ws.Range("F5").Value = 0.00166198
$ws.Range("G5").Value = 0.03787804
$ws.Range("H5").Value = 1.78094815
$ws.Range("I5").Value = 8.199632306998282
$ws.Range("J5").Value = 15.03117028311768
$ws.Range("K5").Value = 0.3118262287250341
$ws.Range("L5").Value = 0.0002498441548279681
$ws.Range("M5").Value = 0.006012809795764376
$ws.Range("N5").Value = 0.5670590400082981
$ws.Range("C6").Value = 13.787
$ws.Range("D6").Value = 55.081
$ws.Range("E6").Value = 3.21703234
$ws.Range("F6").Value = 0.00435844
$ws.Range("G6").Value = 0.01488329
$ws.Range("H6").Value = 0.21439266
$ws.Range("I6").Value = 3.259794051356198
$ws.Range("J6").Value = 13.04640362887979
$ws.Range("K6").Value = 0.656420880163054
$ws.Range("L6").Value = 0.0006457398441681808
$ws.Range("M6").Value = 0.003572530258563019
$ws.Range("N6").Value = 0.1058296747381356
$ws.Range("C7").Value = 24.351
$ws.Range("D7").Value = 90.449
$ws.Range("E7").Value = 1.92385209
$ws.Range("F7").Value = 0.00313757
$ws.Range("G7").Value = 0.01886938
$ws.Range("H7").Value = 0.4741455
$ws.Range("I7").Value = 5.201043150162268
$ws.Range("J7").Value = 15.66961824474151
$ws.Range("K7").Value = 0.3288887290705498
$ws.Range("L7").Value = 0.000655178642931952
$ws.Range("M7").Value = 0.004816035507926994
$ws.Range("N7").Value = 0.2096587577418312
$ws.Range("C8").Value = 9.026
$ws.Range("D8").Value = 54.047
$ws.Range("E8").Value = 3.2465932
$ws.Range("F8").Value = 0.00649477
$ws.Range("G8").Value = 0.00973795
$ws.Range("H8").Value = 0.09139376
$ws.Range("I8").Value = 1.854929854405488
$ws.Range("J8").Value = 11.1118073805868
$ws.Range("K8").Value = 0.6145973977206909
$ws.Range("L8").Value = 0.00153800201472472
$ws.Range("M8").Value = 0.002884756439843269
$ws.Range("N8").Value = 0.04736172435623935
$ws.Range("B9").Value = 0.99995556
$ws.Range("C9").Value = 16.821
$ws.Range("D9").Value = 86.29600000000001
$ws.Range("E9").Value = 2.01509231
$ws.Range("F9").Value = 0.004107620000000001
$ws.Range("G9").Value = 0.01131171
$ws.Range("H9").Value = 0.19793911
$ws.Range("I9").Value = 4.211020267309538
$ws.Range("J9").Value = 14.80106264615671
$ws.Range("K9").Value = 0.3395036868161016
$ws.Range("L9").Value = 0.000663986455896213
$ws.Range("M9").Value = 0.00246834225531821
$ws.Range("N9").Value = 0.09679939926366127
$ws.Range("C10").Value = 6.91
$ws.Range("D10").Value = 55.089
$ws.Range("E10").Value = 3.19726144
$ws.Range("F10").Value = 0.00728126
$ws.Range("G10").Value = 0.00632601
$ws.Range("H10").Value = 0.04580960000000001
$ws.Range("I10").Value = 1.412053077576723
$ws.Range("J10").Value = 11.23250854654686
$ws.Range("K10").Value = 0.662175724919601
$ws.Range("L10").Value = 0.0009135220495067073
$ws.Range("M10").Value = 0.001703840324605493
$ws.Range("N10").Value = 0.02170012876063707
$ws.Range("C11").Value = 12.66
$ws.Range("D11").Value = 80.422
$ws.Range("E11").Value = 2.16211396
$ws.Range("F11").Value = 0.005118800000000001
$ws.Range("G11").Value = 0.00795972
$ws.Range("H11").Value = 0.10568401
$ws.Range("I11").Value = 3.356138206877611
$ws.Range("J11").Value = 13.83876538667476
$ws.Range("K11").Value = 0.3621626749625843
$ws.Range("L11").Value = 0.0008798775736115477
$ws.Range("M11").Value = 0.001979683931861364
$ws.Range("N11").Value = 0.05539014756133772
$ws.Range("C12").Value = 5.608
$ws.Range("D12").Value = 55.843
$ws.Range("E12").Value = 3.180136210000001
$ws.Range("F12").Value = 0.00949055
$ws.Range("G12").Value = 0.005367089999999999
$ws.Range("H12").Value = 0.03192714999999999
$ws.Range("I12").Value = 1.275916911073744
$ws.Range("J12").Value = 12.76625741273958
$ws.Range("K12").Value = 0.7060770387562498
$ws.Range("L12").Value = 0.001716747628488218
$ws.Range("M12").Value = 0.001721453174772097
$ws.Range("N12").Value = 0.01763455825424392
$ws.Range("C13").Value = 10.617
$ws.Range("D13").Value = 75.434
$ws.Range("E13").Value = 2.31641947
$ws.Range("F13").Value = 0.00557833
$ws.Range("G13").Value = 0.00581731
$ws.Range("H13").Value = 0.06624899999999999
$ws.Range("I13").Value = 3.315615222595092
$ws.Range("J13").Value = 14.03512125450395
$ws.Range("K13").Value = 0.4230620195136572
$ws.Range("L13").Value = 0.0009005373276986239
$ws.Range("M13").Value = 0.001655029542541971
$ws.Range("N13").Value = 0.04440052480999975